$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$headers = @("Rep Firm Name", "Brand Carried", "Product Covered", "Product Space")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows: Rep Firm, Brand Carried, Product Covered, Product Space
$data = @(
    @("Example Rep Firm", "ABB", "Variable Frequency Drives", "Flow Control"),
    @("Example Rep Firm", "Old Castle", "One Lift Package Pump Station", "Flow Control"),
    @("Example Rep Firm", "Netzsch", "Progressive Cavity Pumps", "Flow Control"),
    @("Example Rep Firm", "Netzsch", "Rotary Lobe Pumps", "Flow Control"),
    @("Example Rep Firm", "USCP", "Steel Reinforced Polymer Concrete Manholes", "Water Treatment"),
    @("Example Rep Firm", "USCP", "Microtunnel Pipe", "Water Treatment"),
    @("Example Rep Firm", "USCP", "Industrial Pipe Structures", "Water Treatment"),
    @("Example Rep Firm", "Flygt", "Submersible Pumps", "Flow Control"),
    @("Example Rep Firm", "Flygt", "Mixers", "Aeration"),
    @("Example Rep Firm", "Flygt", "Controls", "Flow Control"),
    @("Example Rep Firm", "Flygt", "Check Valves", "Flow Control"),
    @("Example Rep Firm", "Flygt", "Mix-Flush Valves", "Flow Control"),
    @("Example Rep Firm", "ITT Gould Pumps", "Standard Cast Iron Pumps", "Flow Control"),
    @("Example Rep Firm", "ITT Gould Pumps", "Bronze Pumps", "Flow Control"),
    @("Example Rep Firm", "ITT Gould Pumps", "End Suction Pumps", "Flow Control"),
    @("Example Rep Firm", "ITT Gould Pumps", "Vertical Turbine Pumps", "Flow Control"),
    @("Example Rep Firm", "ITT Gould Pumps", "Split Case Pumps", "Flow Control"),
    @("Example Rep Firm", "E/One", "Packaged Low Pressure Sewer Systems", "Flow Control"),
    @("Example Rep Firm", "E/One", "Collection Basin", "Flow Control"),
    @("Example Rep Firm", "E/One", "Grinder Pumps", "Flow Control"),
    @("Example Rep Firm", "Lakeside Equipment", "Equipment for all stages of wastewater treatment", "Wastewater Treatment"),
    @("Example Rep Firm", "Next Turbo", "Geared Turbo Compressors", "Aeration"),
    @("Example Rep Firm", "USF Fabrication", "Aluminum Access Hatches", "Water Treatment"),
    @("Example Rep Firm", "USF Fabrication", "Fall Through Safety Grate System", "Water Treatment")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Apply header style (bold, border, centered/top aligned) from A1 to B1:D1
# (A1 already carries the original header style from before the edit)
$ws.Range("A1").Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove the old row 26 leftover content (clear any extra content beyond D25)
$ws.Range("A26:D26").Clear()
